# Update existing rows 2-7 and add new rows 8-10 for the ECs/FAPs/sCs x ECs/FAPs/sCs
# Dnajb11-Prtg ligand-receptor cross table (per commit: "Natmi following Dr Hou advice").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{ Row=2; A="ECs"; B="Dnajb11"; C="Prtg"; D="ECs"; E=3; F=1; G=9.086337666666665; H=27.259013; I=0.2905360353722418; J=0.2905360353722418; K=2; L=0.6666666666666666; M=0.1787656666666667; N=0.536297; O=0.1283890850233522; P=0.1283890850233522; Q=1.624325210540111; R=14.618926894861; S=0.03730165574775441; T=0.03730165574775442 },
    @{ Row=3; A="ECs"; B="Dnajb11"; C="Prtg"; D="FAPs"; E=3; F=1; G=9.086337666666665; H=27.259013; I=0.2905360353722418; J=0.2905360353722418; K=3; L=1; M=0.7373906666666666; N=2.212172; O=0.5295922576376132; P=0.5295922576376132; Q=6.700180589581777; R=60.30162530623598; S=0.1538656348978669; T=0.153865634897867 },
    @{ Row=4; A="ECs"; B="Dnajb11"; C="Prtg"; D="sCs"; E=3; F=1; G=9.086337666666665; H=27.259013; I=0.2905360353722418; J=0.2905360353722418; K=3; L=1; M=0.476218; N=1.428654; O=0.3420186573390345; P=0.3420186573390345; Q=4.327077550944666; R=38.94369795850199; S=0.09936874472662037; T=0.09936874472662038 },
    @{ Row=5; A="FAPs"; B="Dnajb11"; C="Prtg"; D="ECs"; E=3; F=1; G=13.64094733333333; H=40.922842; I=0.4361698741933416; J=0.4361698741933417; K=2; L=0.6666666666666666; M=0.1787656666666667; N=0.536297; O=0.1283890850233522; P=0.1283890850233522; Q=2.438533044008223; R=21.946797396074; S=0.05599945106243377; T=0.05599945106243378 },
    @{ Row=6; A="FAPs"; B="Dnajb11"; C="Prtg"; D="FAPs"; E=3; F=1; G=13.64094733333333; H=40.922842; I=0.4361698741933416; J=0.4361698741933417; K=3; L=1; M=0.7373906666666666; N=2.212172; O=0.5295922576376132; P=0.5295922576376132; Q=10.05870724809156; R=90.528365232824; S=0.2309921883875655; T=0.2309921883875656 },
    @{ Row=7; A="FAPs"; B="Dnajb11"; C="Prtg"; D="sCs"; E=3; F=1; G=13.64094733333333; H=40.922842; I=0.4361698741933416; J=0.4361698741933417; K=3; L=1; M=0.476218; N=1.428654; O=0.3420186573390345; P=0.3420186573390345; Q=6.496064657185333; R=58.464581914668; S=0.1491782347433423; T=0.1491782347433423 },
    @{ Row=8; A="sCs"; B="Dnajb11"; C="Prtg"; D="ECs"; E=3; F=1; G=8.547106333333334; H=25.641319; I=0.2732940904344165; J=0.2732940904344166; K=2; L=0.6666666666666666; M=0.1787656666666667; N=0.536297; O=0.1283890850233522; P=0.1283890850233522; Q=1.527929161749222; R=13.751362455743; S=0.03508797821316401; T=0.03508797821316402 },
    @{ Row=9; A="sCs"; B="Dnajb11"; C="Prtg"; D="FAPs"; E=3; F=1; G=8.547106333333334; H=25.641319; I=0.2732940904344165; J=0.2732940904344166; K=3; L=1; M=0.7373906666666666; N=2.212172; O=0.5295922576376132; P=0.5295922576376132; Q=6.302556437207556; R=56.723007934868; S=0.1447344343521807; T=0.1447344343521807 },
    @{ Row=10; A="sCs"; B="Dnajb11"; C="Prtg"; D="sCs"; E=3; F=1; G=8.547106333333334; H=25.641319; I=0.2732940904344165; J=0.2732940904344166; K=3; L=1; M=0.476218; N=1.428654; O=0.3420186573390345; P=0.3420186573390345; Q=4.070285883847333; R=36.632572954626; S=0.09347167786907183; T=0.09347167786907185 }
)

$colMap = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20
}

foreach ($entry in $rowsData) {
    $r = $entry.Row
    foreach ($col in $colMap.Keys) {
        $ws.Cells.Item($r, $colMap[$col]).Value = $entry[$col]
    }
}